$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Populate the two new rows (26 and 27). Columns that repeat text already
#    used elsewhere in the sheet are copied (values only) from the existing
#    cell so that the shared-string entry is reused exactly (including any
#    special characters); brand-new text is typed directly.
# ---------------------------------------------------------------------------

# --- Row 26 -----------------------------------------------------------------
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A26").PasteSpecial(-4163) | Out-Null        # xlPasteValues -> "Ministerio de Trabajo y Desarrollo Social"
$ws.Range("B26").Value = 25
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C26").PasteSpecial(-4163) | Out-Null         # xlPasteValues -> "Trabajo"
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4163) | Out-Null         # xlPasteValues -> Mitradel description
$ws.Range("E26").Value = "https://www.mitradel.gob.pa/mitradel-extiende-vigencia-de-los-permisos-de-trabajo/"
$ws.Range("F26").Value = "Considerando que el país se encuentra sufriendo una crisis sanitaria y se mantienen las instrucciones respecto al distanciamiento social y las aglomeraciones, el Ministerio de Trabajo y Desarrollo Laboral (Mitradel), ordena por medio de la Resolución No. DM-163 del 18 de junio de 2020 extender la vigencia de permisos de trabajo, la reapertura de términos dentro de la Dirección Nacional de Empleo y las Direcciones Regionales de Trabajo."
$ws.Range("G4").Copy() | Out-Null
$ws.Range("G26").PasteSpecial(-4163) | Out-Null         # xlPasteValues -> "https://www.mitradel.gob.pa"
$ws.Range("H26").Value = "20-06-2020"
$ws.Range("I26").Value = "19-06-2020"
$ws.Range("J2").Copy() | Out-Null
$ws.Range("J26").PasteSpecial(-4163) | Out-Null         # xlPasteValues -> "Panamá"
$ws.Range("K2").Copy() | Out-Null
$ws.Range("K26").PasteSpecial(-4163) | Out-Null         # xlPasteValues -> "Ministerial"

# --- Row 27 -----------------------------------------------------------------
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A27").PasteSpecial(-4163) | Out-Null
$ws.Range("B27").Value = 26
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4163) | Out-Null
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163) | Out-Null
$ws.Range("E27").Value = "https://www.mitradel.gob.pa/servicios-en-el-mitradel-seran-reactivados-a-partir-del-lunes-22-de-junio/"
$ws.Range("F27").Value = "Las actividades inherentes a las Juntas de Conciliación y Decisión, Dirección de Trabajo y la Dirección de Empleo del Ministerio de Trabajo y Desarrollo Laboral (Mitradel), serán reactivadas a partir del lunes 22 de junio de 2020, atendiendo a las recomendaciones del contenido del documento denominado “El Protocolo para conservar la higiene y salud en el ámbito laboral para la prevención ante el COVID-19” de mantener el mínimo riesgo de contaminación y proteger el bienestar de los usuarios y funcionarios durante el desempeño de sus labores."
$ws.Range("G4").Copy() | Out-Null
$ws.Range("G27").PasteSpecial(-4163) | Out-Null
$ws.Range("H27").Value = "20-06-2020"
$ws.Range("I27").Value = "19-06-2020"
$ws.Range("J2").Copy() | Out-Null
$ws.Range("J27").PasteSpecial(-4163) | Out-Null
$ws.Range("K2").Copy() | Out-Null
$ws.Range("K27").PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Hyperlinks for the new "Descarga Link" (E) and "Sitio Web" (G) cells.
#    (Adding a hyperlink resets the cell style to Excel's built-in
#    "Hyperlink" style, so the table formatting below is (re)applied
#    afterwards.)
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("G26"), "https://www.mitradel.gob.pa", [System.Type]::Missing, [System.Type]::Missing, [System.Type]::Missing) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E26"), "https://www.mitradel.gob.pa/mitradel-extiende-vigencia-de-los-permisos-de-trabajo/", [System.Type]::Missing, [System.Type]::Missing, [System.Type]::Missing) | Out-Null
$ws.Hyperlinks.Add($ws.Range("G27"), "https://www.mitradel.gob.pa", [System.Type]::Missing, [System.Type]::Missing, [System.Type]::Missing) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E27"), "https://www.mitradel.gob.pa/servicios-en-el-mitradel-seran-reactivados-a-partir-del-lunes-22-de-junio/", [System.Type]::Missing, [System.Type]::Missing, [System.Type]::Missing) | Out-Null

# ---------------------------------------------------------------------------
# 3. Apply the formatting (fonts/fills/borders/number formats) of row 21 -
#    which already carries the "last row of table" border treatment that the
#    two brand new rows (26 and 27) must use - onto the two new rows. This is
#    done last so that it overrides the generic "Hyperlink" style Excel
#    applied automatically in step 2.
# ---------------------------------------------------------------------------
$ws.Range("A21:K21").Copy() | Out-Null
$ws.Range("A26:K26").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A21:K21").Copy() | Out-Null
$ws.Range("A27:K27").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Both new data rows wrap long text and need the same 120pt row height as the
# rest of the table's data rows.
$ws.Rows.Item(26).RowHeight = 120
$ws.Rows.Item(27).RowHeight = 120

# ---------------------------------------------------------------------------
# 4. Extend the "Trabajo_CL32" table and its autofilter to cover the two new
#    rows.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:K27")) | Out-Null

# ---------------------------------------------------------------------------
# 5. Extend the Categoria data-validation down through the new rows.
# ---------------------------------------------------------------------------
$validatedRange = $ws.Range("C2:C27")
$validatedRange.Validation.Delete() | Out-Null
$validatedRange.Validation.Add(0, 1, [System.Type]::Missing, [System.Type]::Missing) | Out-Null
$validatedRange.Validation.IgnoreBlank = $true
$validatedRange.Validation.InCellDropdown = $true
$validatedRange.Validation.InputTitle = "Categoria"
$validatedRange.Validation.InputMessage = "Selecciona una categoría de la lista"
$validatedRange.Validation.ErrorTitle = "Entrada no válida"
$validatedRange.Validation.ErrorMessage = "Selecciona una categoría de la lista"
$validatedRange.Validation.ShowInput = $true
$validatedRange.Validation.ShowError = $true

# ---------------------------------------------------------------------------
# 6. Update the visible selection to rest on the newly-added last cell.
# ---------------------------------------------------------------------------
$ws.Range("I27").Select() | Out-Null
